# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

function Set-ForcedText($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-PlainText "D2" "94.738.15"
Set-PlainText "E2" "  +2.58%  "
Set-PlainText "D3" "3.116.63"
Set-PlainText "E3" "  +0.37%  "
Set-ForcedText "D5" "238.27"
Set-PlainText "E5" "  -1.19%  "
Set-ForcedText "D6" "614.04"
Set-PlainText "E6" "  -0.22%  "
Set-PlainText "E7" "  +2.48%  "
Set-ForcedText "D8" "0.392"
Set-PlainText "E8" "  -0.11%  "
Set-ForcedText "D9" "0.999"
Set-PlainText "E9" "  -0.01%  "
Set-ForcedText "D10" "0.837"
Set-PlainText "E10" "  +14.72%  "
Set-PlainText "D11" "3.106.58"
Set-PlainText "E11" "  +0.22%  "
Set-PlainText "E12" "  -2.05%  "
Set-PlainText "E13" "  -2.20%  "
Set-PlainText "D14" "94.292.80"
Set-PlainText "E14" "  +2.62%  "
Set-ForcedText "D15" "34.73"
Set-PlainText "E15" "  +0.99%  "
Set-ForcedText "D16" "5.31"
Set-PlainText "E16" "  -3.59%  "
Set-PlainText "D17" "3.683.99"
Set-PlainText "E17" "  +0.19%  "
Set-PlainText "D18" "3.103.14"
Set-PlainText "E18" "  +0.45%  "
Set-ForcedText "D19" "3.69"
Set-PlainText "E19" "  +0.79%  "
Set-ForcedText "D20" "14.95"
Set-PlainText "E20" "  +1.41%  "
Set-PlainText "E21" "  +1.73%  "
Set-ForcedText "D22" "450.38"
Set-PlainText "E22" "  +0.61%  "
Set-PlainText "E23" "  -1.41%  "
Set-ForcedText "D24" "8.98"
Set-PlainText "E24" "  -4.05%  "
Set-ForcedText "D25" "8.32"
Set-PlainText "E25" "  +5.65%  "
Set-PlainText "E26" "  +0.20%  "
Set-ForcedText "D27" "85.97"
Set-PlainText "E27" "  -1.30%  "
Set-ForcedText "D28" "12.13"
Set-PlainText "E28" "  +3.35%  "
Set-PlainText "D29" "3.286.14"
Set-PlainText "E29" "  +0.60%  "
Set-PlainText "E30" "  -0.12%  "
Set-PlainText "E31" "  +8.37%  "
Set-ForcedText "D32" "0.182"
Set-PlainText "E32" "  +8.77%  "
Set-PlainText "E33" "  -9.78%  "
Set-ForcedText "D34" "9.32"
Set-PlainText "E34" "  +1.10%  "
Set-PlainText "E35" "  +0.26%  "
Set-ForcedText "D36" "7.90"
Set-PlainText "E36" "  -1.65%  "
Set-ForcedText "D37" "0.160"
Set-PlainText "E37" "  +0.08%  "
Set-ForcedText "D38" "26.00"
Set-PlainText "E38" "  -0.78%  "
Set-PlainText "B39" "PolygonEcosystemToken"
Set-PlainText "C39" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-ForcedText "D39" "0.457"
Set-PlainText "E39" "  +5.39%  "
Set-PlainText "B40" "PancakeSwap"
Set-PlainText "C40" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-ForcedText "D40" "1.91"
Set-PlainText "E40" "  -1.02%  "
Set-PlainText "E41" "  +4.93%  "
Set-PlainText "E42" "  -1.09%  "
Set-ForcedText "D43" "472.79"
Set-PlainText "E43" "  -1.59%  "
Set-ForcedText "D44" "3.66"
Set-PlainText "E44" "  -13.99%  "
Set-ForcedText "D45" "3.22"
Set-PlainText "E45" "  -6.68%  "
Set-ForcedText "D47" "160.39"
Set-PlainText "E47" "  +0.60%  "
Set-ForcedText "D48" "0.692"
Set-PlainText "E48" "  -0.29%  "
Set-ForcedText "D49" "1.85"
Set-PlainText "E49" "  -2.31%  "
Set-ForcedText "D50" "4.42"
Set-PlainText "E50" "  +0.35%  "
Set-PlainText "B51" "OKB"
Set-PlainText "C51" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-ForcedText "D51" "43.91"
Set-PlainText "E51" "  -0.53%  "
